# Stationarity test with two similar subtasks: changes made inside
# stationarity_test/stationarity test with subtask_1
#
# - Rename Sheet1 -> Sheet_1
# - Re-number the "Subtask" (column I) groups on Sheet_1: the block that
#   used to be subtask 2 (rows 29-38) becomes subtask 1, and the blocks
#   that used to be subtasks 1, 3, 4 and 5 (rows 2-24, 49-56, 57-66,
#   67-74) all collapse to subtask 0.
# - Inside the new subtask_1 block, move the "countries in general" note
#   from J32 to J31, and add a new note in J30 describing a second,
#   similar subtask about checkin density.
# - Update the active cell / selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the worksheet -------------------------------------------------
$ws.Name = "Sheet_1"

# --- Column I (Subtask) renumbering ---------------------------------------
# rows 2-24 : 1 -> 0
$ws.Range("I2:I24").Value = 0

# rows 29-38 : 2 -> 1
$ws.Range("I29:I38").Value = 1

# rows 49-56 : 3 -> 0
$ws.Range("I49:I56").Value = 0

# rows 57-66 : 4 -> 0
$ws.Range("I57:I66").Value = 0

# rows 67-74 : 5 -> 0
$ws.Range("I67:I74").Value = 0

# --- Column J notes inside the (new) subtask_1 block -----------------------
# J32 ("User begins looking at countries in general") moves up to J31
$ws.Cells.Item(32, 10).ClearContents()
$ws.Cells.Item(31, 10).Value = "User begins looking at countries in general"

# J30 gets a brand-new note describing the second subtask
$ws.Cells.Item(30, 10).Value = "He is trying to figure out which places have denser checkins in the world?"

# --- Sheet view / selection -------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A74").Select() | Out-Null
